{"js": "// Remove the \"Mode -\" definition paragraph (duplicate/erroneous entry,\n// per commit message \"update double file error for prep 3 and prep 4\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text.trim() === \"Mode -\") {\n    paragraph.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Mode -\" definition paragraph (duplicate/erroneous entry,\n# per commit message \"update double file error for prep 3 and prep 4\").\n$d = $word.ActiveDocument\n$paragraphs = $d.Paragraphs\n\nfor ($i = $paragraphs.Count; $i -ge 1; $i--) {\n    $paragraph = $paragraphs.Item($i)\n    if ($paragraph.Range.Text.Trim() -eq \"Mode -\") {\n        $paragraph.Range.Delete()\n    }\n}\n"}
